$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-6 from 45174 (2023-09-05)
# to 45175 (2023-09-06), as serial date numbers.
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45175
}
